$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New rows (4 & 5) for Spanish ("spa") language, mirroring the existing
# "eng" rows 2 & 3.
$ws.Range("A4").Value = "spa"
$ws.Range("B4").Value = 2000003
$ws.Range("C4").Value = "BLZ"
$ws.Range("E4").Value = "Public Holiday"
$ws.Range("F4").Value = "Public Holiday"

$ws.Range("A5").Value = "spa"
$ws.Range("B5").Value = 2000004
$ws.Range("C5").Value = "BZ"
$ws.Range("E5").Value = "Public Holiday"
$ws.Range("F5").Value = "Public Holiday"

# Reuse the same date-number-format style that D2:D3 already use (built-in
# numFmtId 14) by copying it across instead of re-typing the format string,
# which would otherwise register a brand-new custom numFmt.
$ws.Range("D2").Copy()
$ws.Range("D4:D5").PasteSpecial(-4122)
$ws.Range("D4").Value = 45463
$ws.Range("D5").Value = 45556

# Column G ("is_active") holds the literal text "TRUE" (not a boolean) in
# the existing rows. Assigning the string "TRUE" directly auto-converts to
# a real boolean, so copy the value over from an existing "TRUE" cell
# instead (xlPasteValues keeps the original text type and leaves styles
# untouched).
$ws.Range("G2").Copy()
$ws.Range("G4").PasteSpecial(-4163)
$ws.Range("G5").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("A7").Select()
